$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new "Price" values are plain decimal-looking strings (e.g. "1.010").
# Assigning such a string straight to .Value lets Excel auto-convert it into a
# floating point number (losing the trailing zero / exact text form), so for those
# cells we temporarily switch the cell to Text format, write the literal string,
# then restore the cells original style so no formatting changes are left behind.

$ws.Range("D2").Value = '28.502.63'
$ws.Range("E2").Value = '  -3.56%  '
$ws.Range("D3").Value = '1.958.71'
$c = $ws.Range("D4")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.010'
$c.Style = $savedStyle
$ws.Range("E4").Value = '  -0.41%  '
$c = $ws.Range("D5")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '321.81'
$c.Style = $savedStyle
$ws.Range("E5").Value = '  -2.42%  '
$c = $ws.Range("D6")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.011'
$c.Style = $savedStyle
$ws.Range("E6").Value = '  -0.22%  '
$c = $ws.Range("D7")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.4767'
$c.Style = $savedStyle
$ws.Range("E7").Value = '  -4.76%  '
$c = $ws.Range("D8")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.4067'
$c.Style = $savedStyle
$ws.Range("E8").Value = '  -3.59%  '
$c = $ws.Range("D9")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '53.53'
$c.Style = $savedStyle
$ws.Range("E9").Value = '  -1.86%  '
$c = $ws.Range("D10")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.08436'
$c.Style = $savedStyle
$ws.Range("E10").Value = '  -6.52%  '
$c = $ws.Range("D11")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.062'
$c.Style = $savedStyle
$ws.Range("E11").Value = '  -4.92%  '
$c = $ws.Range("D12")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '22.18'
$c.Style = $savedStyle
$ws.Range("E12").Value = '  -4.91%  '
$ws.Range("D13").Value = '1.959.38'
$ws.Range("E13").Value = '  -3.90%  '
$c = $ws.Range("D14")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '7.645'
$c.Style = $savedStyle
$ws.Range("E14").Value = '  -4.94%  '
$c = $ws.Range("D15")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '6.196'
$c.Style = $savedStyle
$ws.Range("E15").Value = '  -4.19%  '
$c = $ws.Range("D16")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.012'
$c.Style = $savedStyle
$ws.Range("E16").Value = '  -0.22%  '
$c = $ws.Range("D17")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.00001075'
$c.Style = $savedStyle
$ws.Range("E17").Value = '  -3.58%  '
$c = $ws.Range("D18")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '89.31'
$c.Style = $savedStyle
$ws.Range("E18").Value = '  -5.45%  '
$c = $ws.Range("D19")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.06612'
$c.Style = $savedStyle
$c = $ws.Range("D20")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '18.75'
$c.Style = $savedStyle
$ws.Range("E20").Value = '  -4.54%  '
$c = $ws.Range("D21")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.010'
$c.Style = $savedStyle
$ws.Range("E21").Value = '  -0.24%  '
$c = $ws.Range("D22")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '5.824'
$c.Style = $savedStyle
$ws.Range("E22").Value = '  -2.53%  '
$ws.Range("D23").Value = '28.537.42'
$ws.Range("E23").Value = '  -3.60%  '
$c = $ws.Range("D24")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '11.62'
$c.Style = $savedStyle
$ws.Range("E24").Value = '  -3.08%  '
$c = $ws.Range("D25")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.291'
$c.Style = $savedStyle
$ws.Range("E25").Value = '  -0.56%  '
$ws.Range("D26").Value = '2.218.01'
$ws.Range("E26").Value = '  -2.76%  '
$c = $ws.Range("D27")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '154.45'
$c.Style = $savedStyle
$ws.Range("E27").Value = '  -2.69%  '
$c = $ws.Range("D28")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '20.25'
$c.Style = $savedStyle
$ws.Range("E28").Value = '  -2.34%  '
$c = $ws.Range("D29")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '6.003'
$c.Style = $savedStyle
$ws.Range("E29").Value = '  -5.53%  '
$c = $ws.Range("D30")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.171'
$c.Style = $savedStyle
$ws.Range("E30").Value = '  -5.53%  '
$ws.Range("E31").Value = '  -3.29%  '
$c = $ws.Range("D32")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.9912'
$c.Style = $savedStyle
$ws.Range("E32").Value = '  -6.06%  '
$c = $ws.Range("D33")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.09609'
$c.Style = $savedStyle
$ws.Range("E33").Value = '  -3.47%  '
$c = $ws.Range("D34")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.452'
$c.Style = $savedStyle
$ws.Range("E34").Value = '  -7.20%  '
$c = $ws.Range("D35")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '5.609'
$c.Style = $savedStyle
$ws.Range("E35").Value = '  -3.88%  '
$c = $ws.Range("D36")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.659'
$c.Style = $savedStyle
$ws.Range("E36").Value = '  -3.72%  '
$c = $ws.Range("D37")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.02339'
$c.Style = $savedStyle
$ws.Range("E37").Value = '  -5.22%  '
$c = $ws.Range("D38")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '8.816'
$c.Style = $savedStyle
$ws.Range("E39").Value = '  -2.68%  '
$c = $ws.Range("D40")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.260'
$c.Style = $savedStyle
$ws.Range("E40").Value = '  -3.54%  '
$c = $ws.Range("D41")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.6240'
$c.Style = $savedStyle
$ws.Range("E41").Value = '  -4.73%  '
$c = $ws.Range("D42")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '11.18'
$c.Style = $savedStyle
$ws.Range("E42").Value = '  -4.26%  '
$ws.Range("E43").Value = '  -0.29%  '
$c = $ws.Range("D44")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.1925'
$c.Style = $savedStyle
$ws.Range("E44").Value = '  -5.88%  '
$c = $ws.Range("D45")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.338'
$c.Style = $savedStyle
$ws.Range("E45").Value = '  +2.66%  '
$c = $ws.Range("D46")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.5975'
$c.Style = $savedStyle
$ws.Range("E46").Value = '  -5.86%  '
$c = $ws.Range("D47")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '12.97'
$c.Style = $savedStyle
$ws.Range("E47").Value = '  -3.31%  '
$c = $ws.Range("D48")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.065'
$c.Style = $savedStyle
$ws.Range("E48").Value = '  -5.84%  '
$c = $ws.Range("D49")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.404'
$c.Style = $savedStyle
$ws.Range("E49").Value = '  -3.01%  '
$c = $ws.Range("D50")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.00000000334'
$c.Style = $savedStyle
$ws.Range("E50").Value = '  -1.99%  '
$c = $ws.Range("D51")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.06844'
$c.Style = $savedStyle
$ws.Range("E51").Value = '  -2.04%  '
